$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the first proxy-server data row (row 2).
# Shared-string table must pick up new entries in this order so that
# the resulting indices line up with the target workbook:
#   7 -> "127.0.0.1", 8 -> "ProxyServer_1", 9 -> "000105001"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "ProxyServer_1"
$ws.Range("B2").Value = "000105001"

# C2 reuses the "ProxyServer_1" string and the same text-format style (s="1")
# that A2/B2 already carry, so make sure it keeps a text number format before
# assigning the value.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "ProxyServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 5001

# Update the active selection shown in the sheet view.
$ws.Range("G4").Select()
